# Apply updated values for rows 2-20, columns B:K (col A labels unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.2071908139402366
$arr[0,1] = -2.16421502631588
$arr[0,2] = -1.499949593677062
$arr[0,3] = -0.7304772392529698
$arr[0,4] = -0.8230868241805323
$arr[0,5] = -0.10215989263453
$arr[0,6] = -0.8111118172483847
$arr[0,7] = -0.3190058249259076
$arr[0,8] = 0.1507989223470795
$arr[0,9] = 0.4636241173125255
$ws.Range("B2:K2").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -2.272829558532021
$arr[0,1] = -1.608564125893203
$arr[0,2] = -0.8390917714691106
$arr[0,3] = -0.931701356396673
$arr[0,4] = -0.2107744248506708
$arr[0,5] = -0.9197263494645256
$arr[0,6] = -0.4276203571420484
$arr[0,7] = 0.04218439013093866
$arr[0,8] = 0.3550095850963846
$arr[0,9] = -0.5233929961551953
$ws.Range("B3:K3").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1.400462404925464
$arr[0,1] = -0.6309900505013716
$arr[0,2] = -0.723599635428934
$arr[0,3] = -0.002672703882931809
$arr[0,4] = -0.7116246284967865
$arr[0,5] = -0.2195186361743094
$arr[0,6] = 0.2502861110986777
$arr[0,7] = 0.5631113060641236
$arr[0,8] = -0.3152912751874563
$arr[0,9] = -1.045203290365925
$ws.Range("B4:K4").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.5144918801275712
$arr[0,1] = -0.6071014650551335
$arr[0,2] = 0.1138254664908686
$arr[0,3] = -0.595126458122986
$arr[0,4] = -0.103020465800509
$arr[0,5] = 0.3667842814724781
$arr[0,6] = 0.6796094764379241
$arr[0,7] = -0.1987931048136559
$arr[0,8] = -0.9287051199921248
$arr[0,9] = 0.7334739763975026
$ws.Range("B5:K5").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.5616080510579985
$arr[0,1] = 0.1593188804880037
$arr[0,2] = -0.549633044125851
$arr[0,3] = -0.05752705180337387
$arr[0,4] = 0.4122776954696132
$arr[0,5] = 0.7251028904350592
$arr[0,6] = -0.1532996908165208
$arr[0,7] = -0.8832117059949898
$arr[0,8] = 0.7789673903946376
$arr[0,9] = 0.1597481019993938
$ws.Range("B6:K6").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.1417647591280393
$arr[0,1] = -0.5671871654858154
$arr[0,2] = -0.07508117316333833
$arr[0,3] = 0.3947235741096488
$arr[0,4] = 0.7075487690750948
$arr[0,5] = -0.1708538121764852
$arr[0,6] = -0.9007658273549541
$arr[0,7] = 0.7614132690346732
$arr[0,8] = 0.1421939806394294
$arr[0,9] = 0.379497744259143
$ws.Range("B7:K7").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.4790798465348092
$arr[0,1] = 0.01302614578766791
$arr[0,2] = 0.482830893060655
$arr[0,3] = 0.795656088026101
$arr[0,4] = -0.082746493225479
$arr[0,5] = -0.8126585084039479
$arr[0,6] = 0.8495205879856794
$arr[0,7] = 0.2303012995904356
$arr[0,8] = 0.4676050632101492
$arr[0,9] = 0.3643276933347375
$ws.Range("B8:K8").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.1916007792754515
$arr[0,1] = 0.6614055265484386
$arr[0,2] = 0.9742307215138846
$arr[0,3] = 0.0958281402623046
$arr[0,4] = -0.6340838749161644
$arr[0,5] = 1.028095221473463
$arr[0,6] = 0.4088759330782192
$arr[0,7] = 0.6461796966979327
$arr[0,8] = 0.542902326822521
$arr[0,9] = -0.2486961005069136
$ws.Range("B9:K9").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1.573432754301089
$arr[0,1] = 1.886257949266535
$arr[0,2] = 1.007855368014955
$arr[0,3] = 0.2779433528364856
$arr[0,4] = 1.940122449226113
$arr[0,5] = 1.320903160830869
$arr[0,6] = 1.558206924450583
$arr[0,7] = 1.454929554575171
$arr[0,8] = 0.6633311272457364
$arr[0,9] = 1.273820034913197
$ws.Range("B10:K10").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.9422837133007778
$arr[0,1] = 0.06388113204919779
$arr[0,2] = -0.6660308831292712
$arr[0,3] = 0.9961482132603562
$arr[0,4] = 0.3769289248651124
$arr[0,5] = 0.6142326884848259
$arr[0,6] = 0.5109553186094142
$arr[0,7] = -0.2806431087200204
$arr[0,8] = 0.3298457989474406
$arr[0,9] = 0.1683237681281231
$ws.Range("B11:K11").Value = $arr

$arr = New-Object 'object[,]' 1,9
$arr[0,0] = 0.0678490295623069
$arr[0,1] = -0.6620629856161621
$arr[0,2] = 1.000116110773465
$arr[0,3] = 0.3808968223782215
$arr[0,4] = 0.6182005859979351
$arr[0,5] = 0.5149232161225235
$arr[0,6] = -0.2766752112069113
$arr[0,7] = 0.3338136964605497
$arr[0,8] = 0.1722916656412322
$ws.Range("B12:J12").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.5264228954459207
$arr[0,1] = 1.135756200943707
$arr[0,2] = 0.5165369125484629
$arr[0,3] = 0.7538406761681764
$arr[0,4] = 0.6505633062927647
$arr[0,5] = -0.1410351210366699
$arr[0,6] = 0.4694537866307911
$arr[0,7] = 0.3079317558114735
$ws.Range("B13:I13").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.8949500190880419
$arr[0,1] = 0.2757307306927982
$arr[0,2] = 0.5130344943125118
$arr[0,3] = 0.4097571244371001
$arr[0,4] = -0.3818413028923346
$arr[0,5] = 0.2286476047751264
$arr[0,6] = 0.06712557395580883
$ws.Range("B14:H14").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.2303995154407018
$arr[0,1] = 0.4677032790604154
$arr[0,2] = 0.3644259091850037
$arr[0,3] = -0.427172518144431
$arr[0,4] = 0.18331638952303
$arr[0,5] = 0.02179435870371246
$ws.Range("B15:G15").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0.4008418571243615
$arr[0,1] = 0.2975644872489498
$arr[0,2] = -0.4940339400804848
$arr[0,3] = 0.1164549675869761
$arr[0,4] = -0.04506706323234141
$ws.Range("B16:F16").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 0.2679782848922332
$arr[0,1] = -0.5236201424372015
$arr[0,2] = 0.08686876523025952
$arr[0,3] = -0.07465326558905801
$ws.Range("B17:E17").Value = $arr

$arr = New-Object 'object[,]' 1,3
$arr[0,0] = -0.5417707991668423
$arr[0,1] = 0.06871810850061863
$arr[0,2] = -0.0928039223186989
$ws.Range("B18:D18").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.0506862842519193
$arr[0,1] = -0.1108357465673982
$ws.Range("B19:C19").Value = $arr

$arr = New-Object 'object[,]' 1,1
$arr[0,0] = -0.1624199859130616
$ws.Range("B20:B20").Value = $arr
Write-Output "Applied ifo GDP component analysis preprocessing updates"
